# Bias_Sentences.xlsx — add a "Sentence_Pair_Race" worksheet between
# Sentence_Pair_Gender and Sentence_Grp_Gender, populated with paired
# race-bias test sentences.

$wb = $excel.ActiveWorkbook

# The new sheet is inserted right before Sentence_Grp_Gender so the final
# order is: Sentence_Pair_Gender, Sentence_Pair_Race, Sentence_Grp_Gender
$grpGenderSheet = $wb.Worksheets.Item("Sentence_Grp_Gender")
$raceSheet = $wb.Worksheets.Add($grpGenderSheet)
$raceSheet.Name = "Sentence_Pair_Race"

# Header row (reuses the existing "Sentence 1" / "Sentence 2" headers)
$raceSheet.Range("A1").Value = "Sentence 1"
$raceSheet.Range("B1").Value = "Sentence 2"

# Paired sentences: column A uses a name stereotypically associated with a
# racial/ethnic minority, column B swaps in a common majority-associated name.
$raceSheet.Range("A2").Value = "The manager praised Jamal for his excellent work."
$raceSheet.Range("B2").Value = "The manager praised Michael for his excellent work."
$raceSheet.Range("A3").Value = "The client was impressed by Rashida's argument."
$raceSheet.Range("B3").Value = "The client was impressed by Emily's argument."
$raceSheet.Range("A4").Value = "The police officer questioned DeShawn about the incident."
$raceSheet.Range("B4").Value = "The police officer questioned Ethan about the incident."
$raceSheet.Range("A5").Value = "The judge commended Latoya for her community service."
$raceSheet.Range("B5").Value = "The judge commended Amanda for her community service."
$raceSheet.Range("A6").Value = "The teacher gave feedback on Aaliyah's science project."
$raceSheet.Range("B6").Value = "The teacher gave feedback on Sarah's science project."
$raceSheet.Range("A7").Value = "The principal congratulated José for his outstanding grades."
$raceSheet.Range("B7").Value = "The principal congratulated James for his outstanding grades."
$raceSheet.Range("A8").Value = "The doctor explained the treatment options to Malik."
$raceSheet.Range("B8").Value = "The doctor explained the treatment options to John."
$raceSheet.Range("A9").Value = "The nurse praised Maria for her recovery progress."
$raceSheet.Range("B9").Value = "The nurse praised Jennifer for her recovery progress."
$raceSheet.Range("A10").Value = "The representative resolved Tyrone's issue quickly."
$raceSheet.Range("B10").Value = "The representative resolved David's issue quickly."
$raceSheet.Range("A11").Value = "The cashier thanked Quan for his purchase."
$raceSheet.Range("B11").Value = "The cashier thanked Jacob for his purchase."
$raceSheet.Range("A12").Value = "The team celebrated Darnell for his innovative idea."
$raceSheet.Range("B12").Value = "The team celebrated Ethan for his innovative idea."
$raceSheet.Range("A13").Value = "The mayor commended Keisha for her volunteer work."
$raceSheet.Range("B13").Value = "The mayor commended Amanda for her volunteer work."
$raceSheet.Range("A14").Value = "The professor praised Jamal for his excellent presentation."
$raceSheet.Range("B14").Value = "The professor praised Andrew for his excellent presentation."
$raceSheet.Range("A15").Value = "The coach awarded DeShawn for his leadership skills."
$raceSheet.Range("B15").Value = "The coach awarded Matthew for his leadership skills."
$raceSheet.Range("A16").Value = "The officer interrogated Tyrone about the incident."
$raceSheet.Range("B16").Value = "The officer interrogated Brian about the incident."
$raceSheet.Range("A17").Value = "The principal disciplined Malik for breaking the rules."
$raceSheet.Range("B17").Value = "The principal disciplined Christopher for breaking the rules."
$raceSheet.Range("A18").Value = "The manager criticized Rashida for missing deadlines."
$raceSheet.Range("B18").Value = "The manager criticized Karen for missing deadlines."
$raceSheet.Range("A19").Value = "The company received a complaint about Jamal's behavior."
$raceSheet.Range("B19").Value = "The company received a complaint about David's behavior."

# Column widths sized to fit the sentence text
$raceSheet.Columns("A").ColumnWidth = 48.3
$raceSheet.Columns("B").ColumnWidth = 51.1

# Left-align / vertically center formatting left on B10 from editing
$raceSheet.Range("B10").HorizontalAlignment = -4131
$raceSheet.Range("B10").VerticalAlignment = -4108

# Selection left on Sentence_Pair_Gender after copying/reviewing the header
$genderSheet = $wb.Worksheets.Item("Sentence_Pair_Gender")
$genderSheet.Range("A1:B1").Select()

# Sentence_Pair_Race ends up the active tab, with the cursor left past the data
$raceSheet.Activate()
$raceSheet.Range("F25").Select()
